$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Poll id suffix bump: "PollN-130120" -> "PollN-230120" (rows 2-8, cols A & B) ---
$ws.Range("A2").Value = "Poll1-230120"
$ws.Range("B2").Value = "Poll1-230120"

$ws.Range("A3").Value = "Poll2-230120"
$ws.Range("B3").Value = "Poll2-230120"

$ws.Range("A4").Value = "Poll3-230120"
$ws.Range("B4").Value = "Poll3-230120"

$ws.Range("A5").Value = "Poll4-230120"
$ws.Range("B5").Value = "Poll4-230120"

$ws.Range("A6").Value = "Poll5-230120"
$ws.Range("B6").Value = "Poll5-230120"

$ws.Range("A7").Value = "Poll6-230120"
$ws.Range("B7").Value = "Poll6-230120"

$ws.Range("A8").Value = "Poll7-230120"
$ws.Range("B8").Value = "Poll7-230120"

# --- privateGroup rename: "Public" -> "All Inside Track Members" (F2) ---
$ws.Range("F2").Value = "All Inside Track Members"

# --- Column F widened to fit the longer group names; split off from the shared C:F width ---
$ws.Columns("F").ColumnWidth = 37.4
